$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specification")
$ws.Rows.Item(12).Delete()
